# --- Refresh cryptos price list (GitHub Actions scheduled snapshot) ---
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.250.44"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "2.295.31"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'315.21"
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("D6").Value = "'102.50"
$ws.Range("E6").Value = "  -3.28%  "
$ws.Range("D7").Value = "'0.628"
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.604"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("D10").Value = "'39.59"
$ws.Range("E10").Value = "  -2.05%  "
$ws.Range("D11").Value = "'0.0905"
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("D12").Value = "'8.37"
$ws.Range("E12").Value = "  +1.37%  "
$ws.Range("D13").Value = "'0.106"
$ws.Range("E13").Value = "  +0.66%  "
$ws.Range("D14").Value = "'0.955"
$ws.Range("E14").Value = "  -1.21%  "
$ws.Range("D15").Value = "'15.19"
$ws.Range("E15").Value = "  -1.74%  "
$ws.Range("D16").Value = "2.643.39"
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("D17").Value = "2.274.93"
$ws.Range("E17").Value = "  -1.44%  "
$ws.Range("D18").Value = "42.241.73"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D19").Value = "'7.38"
$ws.Range("E19").Value = "  -2.06%  "
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("D21").Value = "'73.19"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").Value = "'11.75"
$ws.Range("E22").Value = "  +26.13%  "
$ws.Range("E23").Value = "  +3.07%  "
$ws.Range("D24").Value = "'275.10"
$ws.Range("E25").Value = "  -2.46%  "
$ws.Range("E26").Value = "  -0.45%  "
$ws.Range("E27").Value = "  -1.11%  "
$ws.Range("D28").Value = "'2.41"
$ws.Range("E28").Value = "  +5.25%  "
$ws.Range("D29").Value = "'22.70"
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("D30").Value = "'37.32"
$ws.Range("E30").Value = "  +4.72%  "
$ws.Range("D31").Value = "'165.62"
$ws.Range("E31").Value = "  -0.41%  "
$ws.Range("D32").Value = "'0.0871"
$ws.Range("E32").Value = "  -2.17%  "
$ws.Range("D33").Value = "'5.95"
$ws.Range("E33").Value = "  +3.35%  "
$ws.Range("E34").Value = "  +3.07%  "
$ws.Range("D35").Value = "'2.64"
$ws.Range("E35").Value = "  -9.44%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("E38").Value = "  +2.98%  "
$ws.Range("D39").Value = "'3.70"
$ws.Range("E39").Value = "  +1.86%  "
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("D41").Value = "'1.50"
$ws.Range("E41").Value = "  +2.34%  "
$ws.Range("D42").Value = "'70.03"
$ws.Range("E42").Value = "  -1.88%  "
$ws.Range("D43").Value = "'95.10"
$ws.Range("E43").Value = "  -2.95%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").Value = "'0.225"
$ws.Range("E44").Value = "  -1.01%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  -0.20%  "
$ws.Range("D46").Value = "'12.02"
$ws.Range("E46").Value = "  -2.11%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'112.77"
$ws.Range("E47").Value = "  +0.30%  "
$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D48").Value = "'80.01"
$ws.Range("E48").Value = "  +6.12%  "
$ws.Range("D49").Value = "'8.98"
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("D50").Value = "'5.24"
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("D51").Value = "1.590.57"
$ws.Range("E51").Value = "  +1.96%  "

# Column D values above that are plain decimals (e.g. "102.50") were
# written with a leading apostrophe so Excel keeps the exact text
# (matching the source price format, incl. trailing zeros) instead of
# auto-converting them to a Double. Re-apply the plain number format
# from an untouched price cell so the quote-prefix style marker does
# not linger on those cells.
$protectedCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D19", "D21", "D22", "D24", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D39", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
$protectedRange = $ws.Range($protectedCells[0])
for ($i = 1; $i -lt $protectedCells.Length; $i++) {
    $protectedRange = $excel.Union($protectedRange, $ws.Range($protectedCells[$i]))
}
$ws.Range("D2").Copy()
$protectedRange.PasteSpecial(-4122)
$excel.CutCopyMode = 0
